# Update cryptocurrency price/volume data per Aug 26 2024 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.687.25"
$ws.Range("E2").Value = "  -0.29%  "

# Row 3
$ws.Range("D3").Value = "2.728.20"
$ws.Range("E3").Value = "  -0.78%  "

# Row 4
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.28"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.37%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "160.02"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.10%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.593"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.67%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.108"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.62%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.168"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.76%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.58"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.83%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.373"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.19%  "

# Row 13
$ws.Range("D13").Value = "3.210.64"
$ws.Range("E13").Value = "  -0.73%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.62"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.52%  "

# Row 15
$ws.Range("D15").Value = "63.531.23"
$ws.Range("E15").Value = "  -0.01%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000148"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.02%  "

# Row 17
$ws.Range("D17").Value = "2.729.48"
$ws.Range("E17").Value = "  -0.84%  "

# Row 18
$ws.Range("E18").Value = "  +0.68%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.69"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.66%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "353.01"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.70%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.52"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.92%  "

# Row 22
$ws.Range("E22").Value = "  +0.04%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.515"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.65%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.86"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.18%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.169"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.12%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.02%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.24"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.55%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0896"
$ws.Range("E28").Value = "  -1.46%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.96"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.31%  "

# Row 30
$ws.Range("E30").Value = "  +6.01%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.13"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.42%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "166.66"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.49%  "

# Row 33
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.49"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.15%  "

# Row 34
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.87"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.85%  "

# Row 35
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.93"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.22%  "

# Row 36
$ws.Range("B36").Value = "USDe"
$ws.Range("C36").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.01%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.79"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.16%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "344.14"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.88%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.954"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.29%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.25"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.04%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.06"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.72%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.45"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.38%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.42"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.79%  "

# Row 44
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0578"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.91%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.70"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.89%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.626"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.12%  "

# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "132.25"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.99%  "

# Row 48
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0248"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.70%  "

# Row 49
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.997"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.18%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0989"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.83%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.08"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.50%  "

